$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.012722240760922432
$ws.Range("C2").Value = 0.0052150641568005085
$ws.Range("D2").Value = 0.0040179891511797905
$ws.Range("E2").Value = 0.0029195602983236313
$ws.Range("F2").Value = 0.000017349158952129073
$ws.Range("I2").Value = 1.2575732469558716
$ws.Range("J2").Value = 0.12723591923713684
$ws.Range("K2").Value = 1.4109612703323364
